$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "lartas_*" header cells (E1:H1) are no longer used - clear their text
# but keep the cells present (they pick up the refreshed style below).
$ws.Range("E1:H1").ClearContents() | Out-Null

# Row 1 grows a touch taller (table header row got re-stamped).
$ws.Rows(1).RowHeight = 15

# Re-stamp the whole header row + the HS-code column with the refreshed
# "Normal" cell style (this is what actually produces the new cellXfs entry
# the table now uses).
$ws.Range("A1:H1").Style = "Normal"
$ws.Range("A2:A10").Style = "Normal"

# Selection moved from the old bottom-of-sheet cell to C7.
$ws.Range("C7").Select() | Out-Null
